$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Column widths: the "wide" (12-month / year-end) marker column shifts
# right by one as the data window rolls forward (E, I, M instead of F, J) ---
$ws.Range("E1").EntireColumn.ColumnWidth = 28.15
$ws.Range("F1").EntireColumn.ColumnWidth = 27.15
$ws.Range("I1").EntireColumn.ColumnWidth = 28.15
$ws.Range("J1").EntireColumn.ColumnWidth = 27.15
$ws.Range("M1").EntireColumn.ColumnWidth = 28.15

# --- Row 8: financial-period headers (rolling window shifts left by one
# quarter; a new trailing 12-month period is appended) ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: statement publish-date headers (same rolling shift) ---
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-02-07 (9)"
$ws.Range("F9").Value = "1401-04-26 (4)"
$ws.Range("G9").Value = "1401-08-28 (4)"
$ws.Range("H9").Value = "1401-10-29 (2)"
$ws.Range("I9").Value = "1402-02-13 (9)"
$ws.Range("J9").Value = "1401-04-26 (2)"
$ws.Range("K9").Value = "1401-08-28 (2)"
$ws.Range("L9").Value = "1401-10-29"
$ws.Range("M9").Value = "1402-02-13 (2)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 9292
$ws.Range("E11").Value = 11958
$ws.Range("F11").Value = 3358
$ws.Range("G11").Value = 6689
$ws.Range("H11").Value = 10584
$ws.Range("I11").Value = 15564
$ws.Range("J11").Value = 3812
$ws.Range("K11").Value = 7970
$ws.Range("L11").Value = 12871
$ws.Range("M11").Value = 16912

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) ---
$ws.Range("D12").Value = -5636
$ws.Range("E12").Value = -7121
$ws.Range("F12").Value = -2206
$ws.Range("G12").Value = -3889
$ws.Range("H12").Value = -5706
$ws.Range("I12").Value = -8702
$ws.Range("J12").Value = -1954
$ws.Range("K12").Value = -4177
$ws.Range("L12").Value = -7071
$ws.Range("M12").Value = -9631

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 3656
$ws.Range("E13").Value = 4837
$ws.Range("F13").Value = 1152
$ws.Range("G13").Value = 2799
$ws.Range("H13").Value = 4878
$ws.Range("I13").Value = 6862
$ws.Range("J13").Value = 1858
$ws.Range("K13").Value = 3793
$ws.Range("L13").Value = 5800
$ws.Range("M13").Value = 7281

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -921
$ws.Range("E14").Value = -1012
$ws.Range("F14").Value = -136
$ws.Range("G14").Value = -260
$ws.Range("H14").Value = -394
$ws.Range("I14").Value = -599
$ws.Range("J14").Value = -135
$ws.Range("K14").Value = -309
$ws.Range("L14").Value = -525
$ws.Range("M14").Value = -726

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 2735
$ws.Range("E17").Value = 3825
$ws.Range("F17").Value = 1016
$ws.Range("G17").Value = 2540
$ws.Range("H17").Value = 4484
$ws.Range("I17").Value = 6264
$ws.Range("J17").Value = 1723
$ws.Range("K17").Value = 3485
$ws.Range("L17").Value = 5275
$ws.Range("M17").Value = 6555

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = -7
$ws.Range("G18").Value = -33
$ws.Range("H18").Value = -72
$ws.Range("I18").Value = -99
$ws.Range("J18").Value = -23
$ws.Range("K18").Value = -47
$ws.Range("L18").Value = -71
$ws.Range("M18").Value = -95

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-op income/exp) ---
$ws.Range("D19").Value = 440
$ws.Range("E19").Value = 524
$ws.Range("F19").Value = 238
$ws.Range("G19").Value = 471
$ws.Range("H19").Value = 564
$ws.Range("I19").Value = 691
$ws.Range("J19").Value = 261
$ws.Range("K19").Value = 343
$ws.Range("L19").Value = 744
$ws.Range("M19").Value = 725

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit) ---
$ws.Range("D20").Value = 3175
$ws.Range("E20").Value = 4349
$ws.Range("F20").Value = 1247
$ws.Range("G20").Value = 2978
$ws.Range("H20").Value = 4976
$ws.Range("I20").Value = 6856
$ws.Range("J20").Value = 1961
$ws.Range("K20").Value = 3780
$ws.Range("L20").Value = 5948
$ws.Range("M20").Value = 7184

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -375
$ws.Range("E21").Value = -423
$ws.Range("F21").Value = -125
$ws.Range("G21").Value = -402
$ws.Range("H21").Value = -683
$ws.Range("I21").Value = -789
$ws.Range("J21").Value = -191
$ws.Range("K21").Value = -460
$ws.Range("L21").Value = -750
$ws.Range("M21").Value = -604

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops) ---
$ws.Range("D22").Value = 2800
$ws.Range("E22").Value = 3925
$ws.Range("F22").Value = 1122
$ws.Range("G22").Value = 2576
$ws.Range("H22").Value = 4293
$ws.Range("I22").Value = 6067
$ws.Range("J22").Value = 1771
$ws.Range("K22").Value = 3320
$ws.Range("L22").Value = 5198
$ws.Range("M22").Value = 6581

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 2800
$ws.Range("E24").Value = 3925
$ws.Range("F24").Value = 1122
$ws.Range("G24").Value = 2576
$ws.Range("H24").Value = 4293
$ws.Range("I24").Value = 6067
$ws.Range("J24").Value = 1771
$ws.Range("K24").Value = 3320
$ws.Range("L24").Value = 5198
$ws.Range("M24").Value = 6581

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 568
$ws.Range("E26").Value = 553
$ws.Range("F26").Value = 536
$ws.Range("G26").Value = 504
$ws.Range("H26").Value = 481
$ws.Range("I26").Value = 474
$ws.Range("J26").Value = 681
$ws.Range("K26").Value = 662
$ws.Range("L26").Value = 629
$ws.Range("M26").Value = 567
